$d = $word.ActiveDocument

$pairs = @(
    @("641÷3=213, 2", "243÷4=60, 3"),
    @("130÷3=43, 1", "361÷5=72, 1"),
    @("101÷3=33, 2", "655÷3=218, 1"),
    @("185÷8=23, 1", "971÷5=194, 1"),
    @("522÷5=104, 2", "895÷9=99, 4"),
    @("946÷8=118, 2", "905÷2=452, 1"),
    @("927÷3=309, 0", "170÷2=85, 0"),
    @("350÷9=38, 8", "665÷8=83, 1"),
    @("700÷5=140, 0", "638÷7=91, 1"),
    @("939÷2=469, 1", "764÷3=254, 2"),
    @("660÷8=82, 4", "489÷6=81, 3"),
    @("154÷4=38, 2", "772÷9=85, 7"),
    @("639÷8=79, 7", "834÷4=208, 2"),
    @("670÷9=74, 4", "386÷7=55, 1"),
    @("778÷6=129, 4", "526÷9=58, 4"),
    @("420÷7=60, 0", "981÷3=327, 0"),
    @("807÷4=201, 3", "324÷2=162, 0"),
    @("128÷3=42, 2", "238÷2=119, 0"),
    @("375÷9=41, 6", "686÷4=171, 2"),
    @("137÷8=17, 1", "535÷3=178, 1"),
    @("209÷6=34, 5", "153÷9=17, 0"),
    @("389÷6=64, 5", "661÷6=110, 1"),
    @("283÷4=70, 3", "878÷2=439, 0"),
    @("778÷5=155, 3", "828÷8=103, 4"),
    @("652÷8=81, 4", "996÷2=498, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
